$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.34"
$ws.Range("E2").Value = "'0.12%"
$ws.Range("D3").Value = "'36.95"
$ws.Range("E3").Value = "'3.46%"
$ws.Range("D4").Value = "'5.039"
$ws.Range("E4").Value = "'-1.21%"
$ws.Range("D5").Value = "'0.07898"
$ws.Range("E5").Value = "'0.81%"
$ws.Range("D6").Value = "'2.195"
$ws.Range("E6").Value = "'-3.78%"
$ws.Range("D7").Value = "'8.009"
$ws.Range("E7").Value = "'-0.79%"
$ws.Range("D8").Value = "'4.034"
$ws.Range("E8").Value = "'0.40%"
$ws.Range("D9").Value = "'0.9283"
$ws.Range("E9").Value = "'0.04%"
$ws.Range("D10").Value = "'0.09919"
$ws.Range("E10").Value = "'-0.52%"
$ws.Range("D11").Value = "'0.1881"
$ws.Range("E11").Value = "'3.30%"
$ws.Range("D12").Value = "'0.08672"
$ws.Range("E12").Value = "'0.20%"
$ws.Range("D13").Value = "'0.03603"
$ws.Range("E13").Value = "'6.94%"
$ws.Range("D14").Value = "'0.09949"
$ws.Range("E14").Value = "'0.35%"
$ws.Range("D15").Value = "'0.001483"
$ws.Range("E15").Value = "'-0.47%"
$ws.Range("D16").Value = "'0.005669"
$ws.Range("E16").Value = "'0.63%"
$ws.Range("D17").Value = "'3.457"
$ws.Range("E17").Value = "'-0.80%"
$ws.Range("D18").Value = "'2.483"
$ws.Range("E18").Value = "'18.46%"
$ws.Range("D19").Value = "'0.3435"
$ws.Range("E19").Value = "'0.08%"
$ws.Range("E20").Value = "'0.41%"
$ws.Range("D21").Value = "'4.771"
$ws.Range("E21").Value = "'5.04%"
$ws.Range("D22").Value = "'0.2197"
$ws.Range("E22").Value = "'-1.70%"
$ws.Range("D23").Value = "'0.04596"
$ws.Range("E23").Value = "'-1.42%"
$ws.Range("D24").Value = "'0.005228"
$ws.Range("E24").Value = "'16.55%"
$ws.Range("D25").Value = "'0.001250"
$ws.Range("E25").Value = "'0.91%"
$ws.Range("D26").Value = "'0.0001399"
$ws.Range("E26").Value = "'7.64%"
$ws.Range("D27").Value = "'0.0002714"
$ws.Range("E27").Value = "'0.54%"
$ws.Range("D39").Value = "'0.01828"
$ws.Range("E39").Value = "'3.98%"
$ws.Range("D40").Value = "'0.04770"
$ws.Range("E40").Value = "'1.36%"
$ws.Range("D41").Value = "'0.007982"
$ws.Range("E41").Value = "'1.74%"
$ws.Range("D42").Value = "'0.1414"
$ws.Range("E42").Value = "'-0.15%"
$ws.Range("D43").Value = "'0.007531"
$ws.Range("E43").Value = "'-13.43%"
$ws.Range("D44").Value = "'0.002189"
$ws.Range("E44").Value = "'-4.41%"
$ws.Range("D45").Value = "'0.01043"
$ws.Range("E45").Value = "'13.35%"
$ws.Range("E46").Value = "'4.65%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.09%"
$ws.Range("D48").Value = "'0.0005794"
$ws.Range("E48").Value = "'-0.11%"
$ws.Range("D49").Value = "'36.33"
$ws.Range("E49").Value = "'832.23%"
$ws.Range("D50").Value = "'0.002686"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("D51").Value = "'0.00002098"
$ws.Range("E51").Value = "'-0.09%"
